$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must stay plain text
# (matches original inlineStr typing). Force Text format, set value, then
# restore the Normal style so no stray numFmt/style id is left on the cell.

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = '42.799.39'
$d.Style = "Normal"
$ws.Range("E2").Value = '  +0.23%  '

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = '2.298.11'
$d.Style = "Normal"
$ws.Range("E3").Value = '  -0.41%  '

$d = $ws.Range("D4")
$d.NumberFormat = "@"
$d.Value = '1.00'
$d.Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = '316.88'
$d.Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '

$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = '104.74'
$d.Style = "Normal"
$ws.Range("E6").Value = '  +0.57%  '

$ws.Range("E7").Value = '  -0.66%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -1.25%  '

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = '39.61'
$d.Style = "Normal"
$ws.Range("E10").Value = '  -1.04%  '

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = '0.0910'
$d.Style = "Normal"
$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("E12").Value = '  +1.64%  '

$ws.Range("E13").Value = '  +2.57%  '

$ws.Range("E14").Value = '  +2.89%  '

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = '15.40'
$d.Style = "Normal"
$ws.Range("E15").Value = '  +0.10%  '

$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = '2.647.42'
$d.Style = "Normal"
$ws.Range("E16").Value = '  -0.29%  '

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = '2.321.54'
$d.Style = "Normal"
$ws.Range("E17").Value = '  +0.55%  '

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = '42.736.09'
$d.Style = "Normal"
$ws.Range("E18").Value = '  +0.20%  '

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = '15.18'
$d.Style = "Normal"
$ws.Range("E19").Value = '  +35.98%  '

$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("E21").Value = '  -0.19%  '

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = '74.17'
$d.Style = "Normal"
$ws.Range("E22").Value = '  +1.33%  '

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = '3.55'
$d.Style = "Normal"
$ws.Range("E23").Value = '  -1.15%  '

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = '266.46'
$d.Style = "Normal"
$ws.Range("E24").Value = '  -5.01%  '

$ws.Range("E25").Value = '  -2.49%  '

$ws.Range("E26").Value = '  +0.62%  '

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = '10.93'
$d.Style = "Normal"
$ws.Range("E27").Value = '  +0.46%  '

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = '2.35'
$d.Style = "Normal"
$ws.Range("E28").Value = '  +0.07%  '

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = '6.91'
$d.Style = "Normal"
$ws.Range("E29").Value = '  +16.63%  '

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = '22.61'
$d.Style = "Normal"
$ws.Range("E30").Value = '  -1.47%  '

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = '37.45'
$d.Style = "Normal"
$ws.Range("E31").Value = '  +3.51%  '

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = '166.37'
$d.Style = "Normal"
$ws.Range("E32").Value = '  +0.74%  '

$ws.Range("E33").Value = '  +0.04%  '

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = '0.132'
$d.Style = "Normal"
$ws.Range("E34").Value = '  -3.96%  '

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = '2.62'
$d.Style = "Normal"
$ws.Range("E35").Value = '  +1.04%  '

$ws.Range("E36").Value = '  -3.26%  '

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = '4.57'
$d.Style = "Normal"
$ws.Range("E37").Value = '  -2.07%  '

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = '0.0354'
$d.Style = "Normal"
$ws.Range("E38").Value = '  -6.36%  '

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = '3.75'
$d.Style = "Normal"
$ws.Range("E39").Value = '  -0.04%  '

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = '2.69'
$d.Style = "Normal"
$ws.Range("E40").Value = '  -3.43%  '

$ws.Range("E41").Value = '  +3.76%  '

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = '70.51'
$d.Style = "Normal"
$ws.Range("E42").Value = '  +0.63%  '

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = '0.230'
$d.Style = "Normal"
$ws.Range("E43").Value = '  +0.90%  '

$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = '95.65'
$d.Style = "Normal"
$ws.Range("E44").Value = '  -2.08%  '

$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("E46").Value = '  +1.52%  '

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = '115.01'
$d.Style = "Normal"
$ws.Range("E47").Value = '  +2.44%  '

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = '80.24'
$d.Style = "Normal"
$ws.Range("E48").Value = '  +0.29%  '

$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = '1.710.50'
$d.Style = "Normal"
$ws.Range("E49").Value = '  +6.61%  '

$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = '8.83'
$d.Style = "Normal"
$ws.Range("E50").Value = '  -1.43%  '

$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = '5.11'
$d.Style = "Normal"
$ws.Range("E51").Value = '  -3.64%  '

